$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at positions 191 and 192, pushing the existing
# last row (old row 191) down to row 193.
$ws.Range("A191:A192").EntireRow.Insert()

# Row 191: new "Primera" quality record dated 2023-07-14 (serial 45121)
$ws.Cells.Item(191, 1).Value = 1
$ws.Cells.Item(191, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(191, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(191, 4).Value = 45121
$ws.Cells.Item(191, 5).Value = 15
$ws.Cells.Item(191, 6).Value = "Fruta"
$ws.Cells.Item(191, 7).Value = 100108
$ws.Cells.Item(191, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(191, 9).Value = 100108003
$ws.Cells.Item(191, 10).Value = "Maracuyá"
$ws.Cells.Item(191, 11).Value = "Sin especificar"
$ws.Cells.Item(191, 12).Value = "Primera"
$ws.Cells.Item(191, 13).Value = 220
$ws.Cells.Item(191, 14).Value = 22000
$ws.Cells.Item(191, 15).Value = 23000
$ws.Cells.Item(191, 16).Value = 22545
$ws.Cells.Item(191, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(191, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(191, 19).Value = 1127
$ws.Cells.Item(191, 20).Value = 20

# Row 192: new "Segunda" quality record dated 2023-07-14 (serial 45121)
$ws.Cells.Item(192, 1).Value = 1
$ws.Cells.Item(192, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(192, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(192, 4).Value = 45121
$ws.Cells.Item(192, 5).Value = 15
$ws.Cells.Item(192, 6).Value = "Fruta"
$ws.Cells.Item(192, 7).Value = 100108
$ws.Cells.Item(192, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(192, 9).Value = 100108003
$ws.Cells.Item(192, 10).Value = "Maracuyá"
$ws.Cells.Item(192, 11).Value = "Sin especificar"
$ws.Cells.Item(192, 12).Value = "Segunda"
$ws.Cells.Item(192, 13).Value = 200
$ws.Cells.Item(192, 14).Value = 19000
$ws.Cells.Item(192, 15).Value = 20000
$ws.Cells.Item(192, 16).Value = 19500
$ws.Cells.Item(192, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(192, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(192, 19).Value = 975
$ws.Cells.Item(192, 20).Value = 20

$wb.Save()
